$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - updated TPM-derived values
$ws.Range("M2").Value = 6.460753666666666
$ws.Range("N2").Value = 19.382261
$ws.Range("O2").Value = 0.6826458421750546
$ws.Range("P2").Value = 0.6826458421750545
$ws.Range("Q2").Value = 0.7873418991728888
$ws.Range("R2").Value = 7.086077092556
$ws.Range("S2").Value = 0.6826458421750546
$ws.Range("T2").Value = 0.6826458421750545

# Row 3 - updated TPM-derived values
$ws.Range("O3").Value = 0.2595064250330629
$ws.Range("P3").Value = 0.2595064250330629
$ws.Range("S3").Value = 0.2595064250330629
$ws.Range("T3").Value = 0.2595064250330629

# Row 4 - updated TPM-derived values
$ws.Range("O4").Value = 0.05784773279188247
$ws.Range("P4").Value = 0.05784773279188246
$ws.Range("S4").Value = 0.05784773279188247
$ws.Range("T4").Value = 0.05784773279188246
